# Refresh the cryptocurrency price/volume snapshot (GitHub Actions update).
# Cells in column D that look like plain numbers are forced back to Text
# (NumberFormat "@") before assignment so Excel doesn't silently convert
# them to numeric values (which would drop trailing zeros / change type),
# matching the original inline-string cell type used throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.577.53'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '1.934.34'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.63'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').Value = '  +2.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2913'
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06792'
$ws.Range('E9').Value = '  -0.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '113.55'
$ws.Range('E10').Value = '  +6.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.43'
$ws.Range('E11').Value = '  +4.83%  '
$ws.Range('D12').Value = '1.936.82'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07594'
$ws.Range('E13').Value = '  -1.77%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.490'
$ws.Range('E14').Value = '  +2.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6797'
$ws.Range('E15').Value = '  +0.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '298.62'
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('D17').Value = '30.603.60'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007644'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.190.45'
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9992'
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.551'
$ws.Range('E22').Value = '  -0.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9996'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.514'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.565'
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.42'
$ws.Range('E26').Value = '  +0.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.39'
$ws.Range('E27').Value = '  -2.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.126'
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.431'
$ws.Range('E30').Value = '  +1.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.179'
$ws.Range('E31').Value = '  -0.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.098'
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('E33').Value = '  -1.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7504'
$ws.Range('E34').Value = '  +1.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.148'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02038'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.717'
$ws.Range('E37').Value = '  -0.85%  '
$ws.Range('E38').Value = '  +0.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.031'
$ws.Range('E39').Value = '  -1.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '110.22'
$ws.Range('E40').Value = '  -1.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4463'
$ws.Range('E41').Value = '  -1.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8720'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.830'
$ws.Range('E43').Value = '  -1.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '69.84'
$ws.Range('E45').Value = '  +2.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.321'
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '49.29'
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.301'
$ws.Range('E48').Value = '  -1.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1234'
$ws.Range('E49').Value = '  -2.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.2547'
$ws.Range('E50').Value = '  +2.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.07'
$ws.Range('E51').Value = '  -0.88%  '
